# Refresh cryptos list figures (price + 1h volume change) to match the
# latest GitHub Actions data pull. Row 45/46 also swap coins (Hedera <-> Bittensor)
# since the source ranking reordered those two entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "65.561.38"
$ws.Range("E2").Value = "  -0.49%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.334.75"
$ws.Range("E3").Value = "  -3.63%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5: BNB
$ws.Range("D5").Value = "'576.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.84%  "

# Row 6: Solana
$ws.Range("D6").Value = "'178.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.55%  "

# Row 7: XRP
$ws.Range("E7").Value = "  +3.18%  "

# Row 8: USDC
$ws.Range("E8").Value = "  +0.01%  "

# Row 9: LidoStakedEther
$ws.Range("D9").Value = "3.329.28"
$ws.Range("E9").Value = "  -3.76%  "

# Row 10: Dogecoin
$ws.Range("E10").Value = "  -0.68%  "

# Row 11: Toncoin
$ws.Range("D11").Value = "'6.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.33%  "

# Row 12: Cardano
$ws.Range("D12").Value = "'0.409"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.01%  "

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.908.61"
$ws.Range("E13").Value = "  -3.72%  "

# Row 15: Avalanche
$ws.Range("D15").Value = "'28.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.43%  "

# Row 16: WrappedBTC
$ws.Range("D16").Value = "65.597.72"
$ws.Range("E16").Value = "  -0.57%  "

# Row 17: ShibaInu
$ws.Range("E17").Value = "  -0.77%  "

# Row 18: WrappedEther
$ws.Range("D18").Value = "3.329.68"
$ws.Range("E18").Value = "  -3.86%  "

# Row 19: Polkadot
$ws.Range("D19").Value = "'5.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.03%  "

# Row 20: Chainlink
$ws.Range("D20").Value = "'13.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.68%  "

# Row 21: BitcoinCash
$ws.Range("D21").Value = "'363.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.58%  "

# Row 22: Uniswap
$ws.Range("D22").Value = "'7.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.20%  "

# Row 23: Dai
$ws.Range("E23").Value = "  +0.17%  "

# Row 24: Litecoin
$ws.Range("D24").Value = "'71.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.82%  "

# Row 25: Polygon
$ws.Range("D25").Value = "'0.520"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.12%  "

# Row 26: PEPE
$ws.Range("E26").Value = "  -0.56%  "

# Row 27: InternetComputer(DFINITY)
$ws.Range("D27").Value = "'9.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.56%  "

# Row 28: Kaspa
$ws.Range("E28").Value = "  -0.77%  "

# Row 29: Binance-PegBSC-USD
$ws.Range("E29").Value = "  -0.02%  "

# Row 30: PancakeSwap
$ws.Range("E30").Value = "  -0.72%  "

# Row 31: NEARProtocol
$ws.Range("E31").Value = "  -2.02%  "

# Row 32: USDe
$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.04%  "

# Row 33: EthereumClassic
$ws.Range("D33").Value = "'22.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.78%  "

# Row 34: Aptos
$ws.Range("E34").Value = "  -3.51%  "

# Row 35: Fetch.AI
$ws.Range("E35").Value = "  -4.84%  "

# Row 36: ImmutableX
$ws.Range("E36").Value = "  -1.85%  "

# Row 37: Monero
$ws.Range("D37").Value = "'160.87"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.09%  "

# Row 38: Mantle
$ws.Range("E38").Value = "  -4.87%  "

# Row 39: EnergySwap
$ws.Range("E39").Value = "  -5.88%  "

# Row 40: Stacks
$ws.Range("E40").Value = "  -0.33%  "

# Row 41: dogwifhat
$ws.Range("D41").Value = "'2.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.26%  "

# Row 42: Maker
$ws.Range("D42").Value = "2.712.39"
$ws.Range("E42").Value = "  -2.80%  "

# Row 43: RenderToken
$ws.Range("D43").Value = "'6.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.80%  "

# Row 44: Filecoin
$ws.Range("D44").Value = "'4.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.21%  "

# Row 45: Hedera -> Bittensor
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'336.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.56%  "

# Row 46: Bittensor -> Hedera
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "'0.0668"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.61%  "

# Row 47: OKB
$ws.Range("D47").Value = "'39.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.60%  "

# Row 48: InjectiveProtocol
$ws.Range("D48").Value = "'24.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.04%  "

# Row 49: VeChain
$ws.Range("D49").Value = "'0.0279"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.66%  "

# Row 50: Stellar
$ws.Range("E50").Value = "  +3.04%  "

# Row 51: ONDO
$ws.Range("D51").Value = "'0.967"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.43%  "
